# Update cryptos list with the latest scraped prices/volumes (GitHub Actions refresh).
# Price cells in column D are text-formatted (values like "26.737.94" or "19.40"
# are not valid numbers), so they're entered with a leading apostrophe to force
# text interpretation and then the style is reset to "Normal" so no extra
# number-format style is left attached to the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.737.94"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.22%  "
$ws.Range("D3").Value = "'1.649.20"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.84%  "
$ws.Range("E4").Value = "  +0.31%  "
$ws.Range("D5").Value = "'216.25"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.47%  "
$ws.Range("E6").Value = "  +0.55%  "
$ws.Range("E7").Value = "  +0.21%  "
$ws.Range("E8").Value = "  -0.49%  "
$ws.Range("E9").Value = "  +0.82%  "
$ws.Range("D10").Value = "'19.40"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.22%  "
$ws.Range("D11").Value = "'0.0844"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.43%  "
$ws.Range("D12").Value = "'1.878.04"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.73%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "'1.655.07"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.18%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "'4.23"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.46%  "
$ws.Range("D15").Value = "'0.535"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.58%  "
$ws.Range("D16").Value = "'66.50"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +5.30%  "
$ws.Range("D17").Value = "'26.767.93"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.35%  "
$ws.Range("D18").Value = "'0.0₃0756"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.52%  "
$ws.Range("D19").Value = "'220.30"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.17%  "
$ws.Range("D21").Value = "'4.40"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.27%  "
$ws.Range("D22").Value = "'6.35"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.27%  "
$ws.Range("D23").Value = "'9.58"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.32%  "
$ws.Range("E24").Value = "  +10.49%  "
$ws.Range("D25").Value = "'147.09"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.89%  "
$ws.Range("E26").Value = "  +0.47%  "
$ws.Range("E27").Value = "  -0.01%  "
$ws.Range("E28").Value = "  +2.66%  "
$ws.Range("D29").Value = "'15.91"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.04%  "
$ws.Range("D30").Value = "'0.0519"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.99%  "
$ws.Range("E31").Value = "  +1.20%  "
$ws.Range("D32").Value = "'3.43"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.52%  "
$ws.Range("D33").Value = "'3.07"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.04%  "
$ws.Range("D34").Value = "'1.288.24"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Value = "'1.55"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.94%  "
$ws.Range("E36").Value = "  +7.01%  "
$ws.Range("E37").Value = "  +0.65%  "
$ws.Range("D38").Value = "'0.834"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.99%  "
$ws.Range("D39").Value = "'0.526"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.13%  "
$ws.Range("E40").Value = "  +0.27%  "
$ws.Range("D41").Value = "'0.813"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.34%  "
$ws.Range("E42").Value = "  -1.85%  "
$ws.Range("E43").Value = "  +0.34%  "
$ws.Range("D44").Value = "'1.788.29"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.96%  "
$ws.Range("D45").Value = "'93.80"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.63%  "
$ws.Range("D46").Value = "'59.93"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +9.42%  "
$ws.Range("D47").Value = "'1.61"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.49%  "
$ws.Range("E48").Value = "  +0.81%  "
$ws.Range("D49").Value = "'7.82"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.51%  "
$ws.Range("E50").Value = "  +3.51%  "
$ws.Range("E51").Value = "  -0.46%  "
